$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Estado de Cuenta" detail block lives in rows 16..28 (column E holds the
# "Periodo Mora" label, column F the matching "Valor Mora" amount). The
# commit removes the old period ordering and re-enters the periods/amounts
# in reverse order (newest period first instead of last).
$firstRow = 16
$lastRow = 28

$periods = @()
$amounts = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
    $amounts += $ws.Cells.Item($r, 6).Value2
}

$reversedPeriods = @($periods[($periods.Length - 1)..0])
$reversedAmounts = @($amounts[($amounts.Length - 1)..0])

for ($i = 0; $i -lt $reversedPeriods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $reversedPeriods[$i]
    $ws.Cells.Item($row, 6).Value = $reversedAmounts[$i]
}
